$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "29.072.17"
$ws.Cells.Item(2, 5).Value = "  +0.19%  "

$ws.Cells.Item(3, 4).Value = "1.831.96"
$ws.Cells.Item(3, 5).Value = "  +0.06%  "

$ws.Cells.Item(4, 4).Value = "'1.000"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  +0.22%  "

$ws.Cells.Item(5, 4).Value = "'243.62"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.59%  "

$ws.Cells.Item(6, 4).Value = "'0.6286"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +0.53%  "

$ws.Cells.Item(7, 5).Value = "  +0.19%  "

$ws.Cells.Item(8, 4).Value = "'0.07493"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -0.86%  "

$ws.Cells.Item(9, 4).Value = "'0.2924"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +0.23%  "

$ws.Cells.Item(10, 4).Value = "'23.16"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +2.82%  "

$ws.Cells.Item(11, 4).Value = "'0.07717"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -0.02%  "

$ws.Cells.Item(12, 4).Value = "1.836.84"
$ws.Cells.Item(12, 5).Value = "  +0.38%  "

$ws.Cells.Item(13, 4).Value = "'4.993"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +0.95%  "

$ws.Cells.Item(14, 4).Value = "'0.6683"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +0.72%  "

$ws.Cells.Item(15, 4).Value = "'82.59"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -0.08%  "

$ws.Cells.Item(16, 4).Value = "'0.000009314"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -7.30%  "

$ws.Cells.Item(17, 4).Value = "'6.002"
$ws.Cells.Item(17, 4).Style = "Normal"

$ws.Cells.Item(18, 4).Value = "29.103.37"
$ws.Cells.Item(18, 5).Value = "  +0.29%  "

$ws.Cells.Item(19, 4).Value = "2.080.24"
$ws.Cells.Item(19, 5).Value = "  +0.02%  "

$ws.Cells.Item(20, 4).Value = "'12.58"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +2.12%  "

$ws.Cells.Item(21, 4).Value = "'223.15"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -1.41%  "

$ws.Cells.Item(22, 5).Value = "  +0.40%  "

$ws.Cells.Item(23, 4).Value = "'7.120"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -0.62%  "

$ws.Cells.Item(24, 5).Value = "  +0.21%  "

$ws.Cells.Item(25, 4).Value = "'159.51"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +0.91%  "

$ws.Cells.Item(26, 4).Value = "'0.1399"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +1.98%  "

$ws.Cells.Item(27, 4).Value = "'8.505"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +0.51%  "

$ws.Cells.Item(28, 4).Value = "'17.93"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +0.22%  "

$ws.Cells.Item(29, 4).Value = "'1.497"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +0.66%  "

$ws.Cells.Item(30, 4).Value = "'0.05742"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +10.50%  "

$ws.Cells.Item(31, 4).Value = "'4.153"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +1.68%  "

$ws.Cells.Item(32, 4).Value = "'4.062"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +1.14%  "

$ws.Cells.Item(33, 5).Value = "  +1.46%  "

$ws.Cells.Item(34, 4).Value = "'0.7488"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +1.54%  "

$ws.Cells.Item(35, 4).Value = "'1.847"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +0.17%  "

$ws.Cells.Item(36, 4).Value = "'1.138"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -0.10%  "

$ws.Cells.Item(37, 4).Value = "'2.672"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -0.81%  "

$ws.Cells.Item(38, 4).Value = "'2.762"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +0.15%  "

$ws.Cells.Item(39, 4).Value = "1.220.81"
$ws.Cells.Item(39, 5).Value = "  -1.86%  "

$ws.Cells.Item(40, 4).Value = "'0.01785"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +0.02%  "

$ws.Cells.Item(41, 4).Value = "'6.546"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +3.37%  "

$ws.Cells.Item(42, 4).Value = "'0.8931"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -0.43%  "

$ws.Cells.Item(43, 4).Value = "'1.002"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +0.26%  "

$ws.Cells.Item(44, 4).Value = "'102.17"
$ws.Cells.Item(44, 4).Style = "Normal"

$ws.Cells.Item(45, 4).Value = "'0.00000000127"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +2.99%  "

$ws.Cells.Item(46, 4).Value = "1.984.32"
$ws.Cells.Item(46, 5).Value = "  +0.23%  "

$ws.Cells.Item(47, 4).Value = "'65.68"
$ws.Cells.Item(47, 4).Style = "Normal"

$ws.Cells.Item(48, 4).Value = "'0.07808"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +12.77%  "

$ws.Cells.Item(49, 4).Value = "'0.5089"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -0.29%  "

$ws.Cells.Item(50, 4).Value = "'0.4070"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +0.82%  "

$ws.Cells.Item(51, 4).Value = "'9.009"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +1.54%  "
